$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F24").Value = 148
$ws.Range("G24").Value = 17059.96
$ws.Range("F35").Value = 66
$ws.Range("G35").Value = 3381.18
$ws.Range("F37").Value = 31
$ws.Range("G37").Value = 825.84
$ws.Range("B40").Value = 54256.77
$ws.Range("F42").Value = 75
$ws.Range("G42").Value = 14757.75
$ws.Range("F44").Value = 409
$ws.Range("G44").Value = 14891.69
$ws.Range("F47").Value = 197
$ws.Range("G47").Value = 37999.33
$ws.Range("F57").Value = 126
$ws.Range("G57").Value = 11786.04
$ws.Range("F58").Value = 15
$ws.Range("G58").Value = 885
$ws.Range("F59").Value = 27
$ws.Range("G59").Value = 443.07
$ws.Range("F67").Value = 199
$ws.Range("G67").Value = 51885.27
$ws.Range("B72").Value = 182616.89
$ws.Range("F84").Value = 97
$ws.Range("G84").Value = 3038.04
$ws.Range("B89").Value = 12194.03
$ws.Range("B132").Value = 65258
$ws.Range("F132").Value = 0
$ws.Range("G132").Value = 0
$ws.Range("B133").Value = 64196
$ws.Range("F133").Value = 1
$ws.Range("G133").Value = 32143.58
$ws.Range("F160").Value = 291
$ws.Range("G160").Value = 9701.940000000001
$ws.Range("B161").Value = 34551.63
$ws.Range("F181").Value = 27
$ws.Range("G181").Value = 7830.27
$ws.Range("F188").Value = 13
$ws.Range("G188").Value = 1164.02
$ws.Range("B199").Value = 58950.65
$ws.Range("F212").Value = 52
$ws.Range("G212").Value = 3369.6
$ws.Range("B214").Value = 3369.6
$ws.Range("F220").Value = 91
$ws.Range("G220").Value = 9596.860000000001
$ws.Range("B224").Value = 69043.81
$ws.Range("F228").Value = 510
$ws.Range("G228").Value = 9435
$ws.Range("F229").Value = 20
$ws.Range("G229").Value = 428.6
$ws.Range("F233").Value = 25
$ws.Range("G233").Value = 2865
$ws.Range("B235").Value = 18144.02
$ws.Range("F267").Value = 83
$ws.Range("G267").Value = 6806.83
$ws.Range("F268").Value = 67
$ws.Range("G268").Value = 5248.11
$ws.Range("F270").Value = 54
$ws.Range("G270").Value = 1881.36
$ws.Range("F277").Value = 18
$ws.Range("G277").Value = 2062.08
$ws.Range("F287").Value = 1
$ws.Range("G287").Value = 26.89
$ws.Range("B296").Value = 64983
$ws.Range("C296").Value = 'HIM-TOTAL CARE BABY PANTS DIAPERS-M-9S'
$ws.Range("F296").Value = 6
$ws.Range("G296").Value = 514.08
$ws.Range("B297").Value = 66194
$ws.Range("C297").Value = 'HIM-Total Care Baby Pants Diapers-M-9s'
$ws.Range("F297").Value = 22
$ws.Range("G297").Value = 1884.96
$ws.Range("B301").Value = 104140.99
$ws.Range("F355").Value = 134
$ws.Range("G355").Value = 10000.42
$ws.Range("B362").Value = 74847.19
$ws.Range("F368").Value = 21
$ws.Range("G368").Value = 15626.31
$ws.Range("B369").Value = 64491.48
$ws.Range("F377").Value = 62
$ws.Range("G377").Value = 9316.74
$ws.Range("B378").Value = 52659.69
$ws.Range("F393").Value = 388
$ws.Range("G393").Value = 37480.8
$ws.Range("B395").Value = 53630.34
$ws.Range("F402").Value = 118
$ws.Range("G402").Value = 3006.64
$ws.Range("F408").Value = 25
$ws.Range("G408").Value = 857.75
$ws.Range("F409").Value = 64
$ws.Range("G409").Value = 2594.56
$ws.Range("F414").Value = 176
$ws.Range("G414").Value = 2789.6
$ws.Range("F419").Value = 71
$ws.Range("G419").Value = 4087.47
$ws.Range("B423").Value = 159167.78
$ws.Range("F426").Value = 35
$ws.Range("G426").Value = 5562.2
$ws.Range("F429").Value = 52
$ws.Range("G429").Value = 4346.68
$ws.Range("B433").Value = 22859.84
$ws.Range("F436").Value = 209
$ws.Range("G436").Value = 9672.52
$ws.Range("F437").Value = 11
$ws.Range("G437").Value = 295.79
$ws.Range("F439").Value = 96
$ws.Range("G439").Value = 925.4400000000001
$ws.Range("B444").Value = 22129.65
$ws.Range("F460").Value = 59
$ws.Range("G460").Value = 16694.64
$ws.Range("B464").Value = 86567.52
$ws.Range("B485").Value = 53319
$ws.Range("E485").Value = 310.64
$ws.Range("F485").Value = -6
$ws.Range("G485").Value = -1643.52
$ws.Range("B486").Value = 64810
$ws.Range("E486").Value = 291.22
$ws.Range("F486").Value = 0
$ws.Range("G486").Value = 0
$ws.Range("F506").Value = 87
$ws.Range("G506").Value = 8569.5
$ws.Range("B512").Value = 64830
$ws.Range("E512").Value = 34.9
$ws.Range("F512").Value = 83
$ws.Range("G512").Value = 2724.89
$ws.Range("B513").Value = 60022
$ws.Range("E513").Value = 37.22
$ws.Range("F513").Value = -113
$ws.Range("G513").Value = -3709.79
$ws.Range("B514").Value = 40893.9
$ws.Range("F517").Value = 200
$ws.Range("G517").Value = 19974
$ws.Range("F523").Value = 134
$ws.Range("G523").Value = 7942.18
$ws.Range("F527").Value = 79
$ws.Range("G527").Value = 2164.6
$ws.Range("B531").Value = 112231.96
$ws.Range("F533").Value = 27
$ws.Range("G533").Value = 893.97
$ws.Range("F535").Value = 110
$ws.Range("G535").Value = 3642.1
$ws.Range("F537").Value = 189
$ws.Range("G537").Value = 6257.79
$ws.Range("F540").Value = 115
$ws.Range("G540").Value = 5032.4
$ws.Range("B541").Value = 20543.84
$ws.Range("F564").Value = 154
$ws.Range("G564").Value = 18764.9
$ws.Range("B567").Value = 21288.56
$ws.Range("F611").Value = 164
$ws.Range("G611").Value = 21828.4
$ws.Range("B613").Value = 21828.4
$ws.Range("F615").Value = 0
$ws.Range("G615").Value = 0
$ws.Range("F618").Value = 224
$ws.Range("G618").Value = 33691.84
$ws.Range("F621").Value = 95
$ws.Range("G621").Value = 14693.65
$ws.Range("F628").Value = 464
$ws.Range("G628").Value = 47750.24
$ws.Range("B634").Value = 196280.13
$ws.Range("F665").Value = 31
$ws.Range("G665").Value = 1659.74
$ws.Range("B674").Value = 10089.8
$ws.Range("F680").Value = 554
$ws.Range("G680").Value = 90362.94
$ws.Range("B686").Value = 91375.49000000001
$ws.Range("B724").Value = 2462796.97
$ws.Range("B725").Value = 2462796.97
